# "slightly updated some timing exports for SD_CPT (sCPT)"
#
# Pastes refreshed SD_CPT timing numbers into columns P, W, X and Y
# (rows 2-11) of Sheet1 -- column X had held the formula "=P+L" and is
# now overwritten with a plain pasted value like the others. Also fixes
# a typo in one of the footnotes and clears a few stray cells that had
# been accidentally filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Updated timing data, rows 2-11.
# ---------------------------------------------------------------------

$ws.Range("P2").Value = 0.090052
$ws.Range("W2").Value = 0.315695
$ws.Range("X2").Value = 0.252613
$ws.Range("Y2").Value = 0.259575

$ws.Range("P3").Value = 0.413632
$ws.Range("W3").Value = 1.38659
$ws.Range("X3").Value = 0.995095
$ws.Range("Y3").Value = 1.00385

$ws.Range("P4").Value = 0.882699
$ws.Range("W4").Value = 4.08797
$ws.Range("X4").Value = 2.22419
$ws.Range("Y4").Value = 2.23958

$ws.Range("P5").Value = 0.050236
$ws.Range("W5").Value = 0.12792
$ws.Range("X5").Value = 0.097213
$ws.Range("Y5").Value = 0.097602

$ws.Range("P6").Value = 1.36226
$ws.Range("W6").Value = 3.84295
$ws.Range("X6").Value = 2.60126
$ws.Range("Y6").Value = 2.63357

$ws.Range("P7").Value = 0.358605
$ws.Range("W7").Value = 1.57727
$ws.Range("X7").Value = 0.818705
$ws.Range("Y7").Value = 0.824114

$ws.Range("P8").Value = 0.195486
$ws.Range("W8").Value = 0.332171
$ws.Range("X8").Value = 0.321097
$ws.Range("Y8").Value = 0.332032

$ws.Range("P9").Value = 1.14887
$ws.Range("W9").Value = 2.43686
$ws.Range("X9").Value = 2.10142
$ws.Range("Y9").Value = 2.14747

$ws.Range("P10").Value = 0.380455
$ws.Range("W10").Value = 1.23034
$ws.Range("X10").Value = 0.869905
$ws.Range("Y10").Value = 0.891952

$ws.Range("P11").Value = 0.5698
$ws.Range("W11").Value = 1.62714
$ws.Range("X11").Value = 1.20861
$ws.Range("Y11").Value = 1.23839

# ---------------------------------------------------------------------
# 2. Footnote typo fix + clearing stray copy/paste leftovers.
# ---------------------------------------------------------------------

$ws.Range("J13").Value = "This memory report for II contains a further array Sigma x uint64_t for CPT/+ implementation"

$ws.Range("W13").Clear()
$ws.Range("X13").Clear()
$ws.Range("Y13").Clear()

$ws.Range("X16").Clear()
$ws.Range("X17").Clear()
$ws.Range("X18").Clear()
$ws.Range("X19").Clear()
$ws.Range("X20").Clear()
$ws.Range("X21").Clear()
$ws.Range("X22").Clear()
$ws.Range("X23").Clear()
$ws.Range("X24").Clear()
$ws.Range("X25").Clear()

# ---------------------------------------------------------------------
# 3. Updated saved view (scroll position / active cell).
# ---------------------------------------------------------------------

$ws.Activate()
$ws.Range("E1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("J13").Select()

$excel.CalculateFull()
